# Auto-generated Excel COM-interop edit script
# Applies per-cell numeric updates (and a few cell clears) to the
# Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# reflecting refreshed currentAveragePrice-derived market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 13047.056
$ws.Range("I40").Value = 1835.2858
$ws.Range("J40").Value = 20181.818
$ws.Range("K40").Value = 1835.2858
$ws.Range("L40").Value = 20181.818
$ws.Range("M40").Value = -1660.2858
$ws.Range("N40").Value = -20531.818

$ws.Range("H138").Value = 3144.611
$ws.Range("J138").Value = 3978.2188
$ws.Range("L138").Value = 11934.6564
$ws.Range("N138").Value = -22214.6564

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2661.652
$ws.Range("I2").Value = 2382.8333
$ws.Range("J2").Value = 2965.818
$ws.Range("K2").Value = 2382.8333
$ws.Range("L2").Value = 2965.818
$ws.Range("M2").Value = -2269.8333
$ws.Range("N2").Value = -3191.818

$ws.Range("H116").Value = 2661.652
$ws.Range("I116").Value = 2382.8333
$ws.Range("J116").Value = 2965.818
$ws.Range("K116").Value = 2382.8333
$ws.Range("L116").Value = 2965.818
$ws.Range("M116").Value = -88.83329999999978
$ws.Range("N116").Value = -7553.818

$ws.Range("H141").Value = 89422.25
$ws.Range("J141").Value = 89422.25
$ws.Range("L141").Value = 89422.25
$ws.Range("N141").Value = -99782.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2661.652
$ws.Range("I3").Value = 2382.8333
$ws.Range("J3").Value = 2965.818
$ws.Range("K3").Value = 2382.8333
$ws.Range("L3").Value = 2965.818
$ws.Range("M3").Value = -2268.8333
$ws.Range("N3").Value = -3193.818

$ws.Range("H105").Value = 2396.6667
$ws.Range("I105").Value = 2396.6667
$ws.Range("K105").Value = 2396.6667
$ws.Range("M105").Value = -649.6667000000002

$ws.Range("H133").Value = 19999.25
$ws.Range("J133").Value = 19999.25
$ws.Range("L133").Value = 19999.25
$ws.Range("N133").Value = -30119.25

$ws.Range("H139").Value = 73069.5
$ws.Range("J139").Value = 71927
$ws.Range("L139").Value = 71927
$ws.Range("N139").Value = -82207

$ws.Range("H140").Value = 94999
$ws.Range("J140").Value = 94999
$ws.Range("L140").Value = 94999
$ws.Range("N140").Value = -105359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 22500
$ws.Range("J50").Value = 25000
$ws.Range("L50").Value = 25000
$ws.Range("N50").Value = -26250

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""  # was -41472

$ws.Range("H59").Value = 133333
$ws.Range("J59").Value = 133333
$ws.Range("L59").Value = 133333
$ws.Range("N59").Value = -135623

$ws.Range("H60").Value = 38000
$ws.Range("I60").Value = 38000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 38000
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -37489
$ws.Range("N60").Value = ""  # was -44022

$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = ""  # was -40696

$ws.Range("H96").Value = 36226.855
$ws.Range("J96").Value = 36226.855
$ws.Range("L96").Value = 36226.855
$ws.Range("N96").Value = -41718.855

$ws.Range("H99").Value = 5560865
$ws.Range("J99").Value = 4900
$ws.Range("L99").Value = 4900
$ws.Range("N99").Value = -7896

$ws.Range("H103").Value = 58749
$ws.Range("I103").Value = 54999
$ws.Range("J103").Value = 62499
$ws.Range("K103").Value = 54999
$ws.Range("L103").Value = 62499
$ws.Range("M103").Value = -53827
$ws.Range("N103").Value = -64843

$ws.Range("H126").Value = 5560865
$ws.Range("J126").Value = 4900
$ws.Range("L126").Value = 14700
$ws.Range("N126").Value = -19640

$ws.Range("H134").Value = 4646.2964
$ws.Range("I134").Value = 4705.44
$ws.Range("K134").Value = 14116.32
$ws.Range("M134").Value = -11581.32

$ws.Range("H137").Value = 69476.60000000001
$ws.Range("J137").Value = 69998.44500000001
$ws.Range("L137").Value = 69998.44500000001
$ws.Range("N137").Value = -80198.44500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 69873.7
$ws.Range("J37").Value = 69873.7
$ws.Range("L37").Value = 209621.1
$ws.Range("N37").Value = -209845.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1948225.6
$ws.Range("I11").Value = 3155172.8
$ws.Range("J11").Value = 258499.6
$ws.Range("K11").Value = 3155172.8
$ws.Range("L11").Value = 258499.6
$ws.Range("M11").Value = -3155033.8
$ws.Range("N11").Value = -258777.6

$ws.Range("H15").Value = 99999
$ws.Range("J15").Value = 99999
$ws.Range("L15").Value = 99999
$ws.Range("N15").Value = -100575

$ws.Range("H81").Value = 99999
$ws.Range("J81").Value = 99999
$ws.Range("L81").Value = 99999
$ws.Range("N81").Value = -101995

$ws.Range("H84").Value = 99999
$ws.Range("J84").Value = 99999
$ws.Range("L84").Value = 299997
$ws.Range("N84").Value = -309981

$ws.Range("H101").Value = 24811.143
$ws.Range("J101").Value = 24811.143
$ws.Range("L101").Value = 24811.143
$ws.Range("N101").Value = -31301.143

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""  # was -14900

$ws.Range("H124").Value = 89998.664
$ws.Range("J124").Value = 89998.664
$ws.Range("L124").Value = 89998.664
$ws.Range("N124").Value = -99818.664

$ws.Range("H137").Value = 88998
$ws.Range("J137").Value = 88998
$ws.Range("L137").Value = 88998
$ws.Range("N137").Value = -99198

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4673.1665
$ws.Range("I46").Value = 5326.0713
$ws.Range("J46").Value = 3759.1
$ws.Range("K46").Value = 5326.0713
$ws.Range("L46").Value = 3759.1
$ws.Range("M46").Value = -5138.0713
$ws.Range("N46").Value = -4135.1

$ws.Range("H138").Value = 99993.5
$ws.Range("J138").Value = 99993.5
$ws.Range("L138").Value = 99993.5
$ws.Range("N138").Value = -110273.5

$ws.Range("H139").Value = 88570.42999999999
$ws.Range("J139").Value = 88570.42999999999
$ws.Range("L139").Value = 88570.42999999999
$ws.Range("N139").Value = -98850.42999999999

$ws.Range("H140").Value = 80610.5
$ws.Range("J140").Value = 74572.78
$ws.Range("L140").Value = 74572.78
$ws.Range("N140").Value = -84932.78

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 89999
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = ""  # was -33685

$ws.Range("H73").Value = 89999
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = ""  # was -32908

$ws.Range("H126").Value = 3999.4
$ws.Range("I126").Value = 3998.5
$ws.Range("K126").Value = 11995.5
$ws.Range("M126").Value = -9525.5
